$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'312.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'5.17%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'44.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'7.47%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.114"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.40%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08007"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'6.04%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.497"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'2.68%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.656"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'3.73%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.084"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'16.90%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1292"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'6.90%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1907"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.71%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09383"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'5.98%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04249"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'7.77%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1039"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.92%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001313"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'2.45%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005762"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-2.64%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D17").Value = "'3.376"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.24%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D19").Value = "'0.3385"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.97%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.067"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.82%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1349"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-4.91%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'4.27%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'3.26%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001273"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.50%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004602"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'15.41%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001338"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'8.69%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02657"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'10.33%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05415"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'3.89%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.005619"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-12.28%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007725"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.75%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1411"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'6.22%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007328"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-3.05%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007869"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'0.35%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3122"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-2.88%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006744"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.63%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000744"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.93%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.05600"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'21.03%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.003966"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-5.64%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002082"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.93%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0001983"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.93%"
$ws.Range("E51").Style = "Normal"
